# Auto-generated script to apply scheduled market-data refresh values
# to the Halicarnassus_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 31
$ws.Range("I11").Value = 31
$ws.Range("K11").Value = 31
$ws.Range("M11").Value = 109
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").Value = $null
$ws.Range("H62").Value = 11201.2
$ws.Range("J62").Value = 11201.2
$ws.Range("L62").Value = 11201.2
$ws.Range("N62").Value = -12449.2
$ws.Range("H65").Value = 11201.2
$ws.Range("J65").Value = 11201.2
$ws.Range("L65").Value = 56006
$ws.Range("N65").Value = -62246
$ws.Range("H69").Value = 7221.4287
$ws.Range("J69").Value = 7221.4287
$ws.Range("L69").Value = 21664.2861
$ws.Range("N69").Value = -23412.2861
$ws.Range("H70").Value = 5089.8887
$ws.Range("I70").Value = 2236
$ws.Range("K70").Value = 6708
$ws.Range("M70").Value = -6438
$ws.Range("H72").Value = 7221.4287
$ws.Range("J72").Value = 7221.4287
$ws.Range("L72").Value = 64992.85830000001
$ws.Range("N72").Value = -73728.85830000001
$ws.Range("H73").Value = 5089.8887
$ws.Range("I73").Value = 2236
$ws.Range("K73").Value = 6708
$ws.Range("M73").Value = -5772
$ws.Range("H74").Value = 4001.5
$ws.Range("I74").Value = 3003
$ws.Range("K74").Value = 3003
$ws.Range("M74").Value = -2067
$ws.Range("H77").Value = 4001.5
$ws.Range("I77").Value = 3003
$ws.Range("K77").Value = 15015
$ws.Range("M77").Value = -10335
$ws.Range("H88").Value = 1032.125
$ws.Range("I88").Value = 958.8333
$ws.Range("J88").Value = 1252
$ws.Range("K88").Value = 958.8333
$ws.Range("L88").Value = 1252
$ws.Range("M88").Value = -552.8333
$ws.Range("N88").Value = -2064
$ws.Range("H91").Value = 1032.125
$ws.Range("I91").Value = 958.8333
$ws.Range("J91").Value = 1252
$ws.Range("K91").Value = 958.8333
$ws.Range("L91").Value = 1252
$ws.Range("M91").Value = 445.1667
$ws.Range("N91").Value = -4060
$ws.Range("H100").Value = 277.33334
$ws.Range("I100").Value = 296.375
$ws.Range("J100").Value = 125
$ws.Range("K100").Value = 296.375
$ws.Range("L100").Value = 125
$ws.Range("M100").Value = 244.625
$ws.Range("N100").Value = -1207
$ws.Range("H137").Value = 2584.6
$ws.Range("I137").Value = 1138.4
$ws.Range("J137").Value = 3066.6667
$ws.Range("K137").Value = 3415.2
$ws.Range("L137").Value = 9200.000100000001
$ws.Range("M137").Value = -865.2000000000003
$ws.Range("N137").Value = -14300.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 8233
$ws.Range("I63").Value = 2350
$ws.Range("K63").Value = 2350
$ws.Range("M63").Value = -1664
$ws.Range("H66").Value = 8233
$ws.Range("I66").Value = 2350
$ws.Range("K66").Value = 11750
$ws.Range("M66").Value = -8318
$ws.Range("H97").Value = 783.5
$ws.Range("I97").Value = 953.6
$ws.Range("K97").Value = 953.6
$ws.Range("M97").Value = -457.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 2699.9092
$ws.Range("I22").Value = 2919.9
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 2919.9
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -2746.9
$ws.Range("N22").Value = -846

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").Value = $null
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").Value = $null
$ws.Range("H74").Value = 48164.11
$ws.Range("J74").Value = 48164.11
$ws.Range("L74").Value = 48164.11
$ws.Range("N74").Value = -49912.11
$ws.Range("H77").Value = 48164.11
$ws.Range("J77").Value = 48164.11
$ws.Range("L77").Value = 144492.33
$ws.Range("N77").Value = -153228.33

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 211414.4
$ws.Range("I4").Value = 400450.2
$ws.Range("J4").Value = 22378.6
$ws.Range("K4").Value = 1201350.6
$ws.Range("L4").Value = 67135.79999999999
$ws.Range("M4").Value = -1201238.6
$ws.Range("N4").Value = -67359.79999999999
$ws.Range("H59").Value = 900
$ws.Range("I59").Value = 900
$ws.Range("K59").Value = 2700
$ws.Range("M59").Value = -2160

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 18890162
$ws.Range("I3").Value = 20864642
$ws.Range("J3").Value = 10005000
$ws.Range("K3").Value = 20864642
$ws.Range("L3").Value = 10005000
$ws.Range("M3").Value = -20864526
$ws.Range("N3").Value = -10005232
$ws.Range("H80").Value = 3745
$ws.Range("I80").Value = 3990
$ws.Range("K80").Value = 3990
$ws.Range("M80").Value = -2992
$ws.Range("H83").Value = 3745
$ws.Range("I83").Value = 3990
$ws.Range("K83").Value = 19950
$ws.Range("M83").Value = -14958
$ws.Range("H122").Value = 3299.4167
$ws.Range("I122").Value = 3036.625
$ws.Range("K122").Value = 9109.875
$ws.Range("M122").Value = -6659.875
$ws.Range("H123").Value = 983333.3
$ws.Range("J123").Value = 950000
$ws.Range("L123").Value = 950000
$ws.Range("N123").Value = -954900
$ws.Range("H126").Value = 2374.75
$ws.Range("J126").Value = 2500
$ws.Range("L126").Value = 7500
$ws.Range("N126").Value = -12440

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4907.727
$ws.Range("I40").Value = 4454.6665
$ws.Range("K40").Value = 4454.6665
$ws.Range("M40").Value = -4318.6665
$ws.Range("H82").Value = 5094.25
$ws.Range("I82").Value = 2892.4
$ws.Range("J82").Value = 6667
$ws.Range("K82").Value = 2892.4
$ws.Range("L82").Value = 6667
$ws.Range("M82").Value = -2531.4
$ws.Range("N82").Value = -7389
$ws.Range("H85").Value = 5094.25
$ws.Range("I85").Value = 2892.4
$ws.Range("J85").Value = 6667
$ws.Range("K85").Value = 2892.4
$ws.Range("L85").Value = 6667
$ws.Range("M85").Value = -1644.4
$ws.Range("N85").Value = -9163
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4081.4
$ws.Range("I122").Value = 2467.5
$ws.Range("K122").Value = 7402.5
$ws.Range("M122").Value = -4952.5
